# Insert two new data rows (159 and 160) into the "Camote" sheet, shifting the
# existing rows 159:184 down to 161:186, and populate the two new rows with
# their data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 159 (Excel copies formatting
# from the row above by default, which matches the desired style for column D).
$ws.Rows("159:160").Insert()

# --- New row 159 ---
$ws.Cells.Item(159, 1).Value2  = 9
$ws.Cells.Item(159, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(159, 3).Value2  = "Metropolitana"
$ws.Cells.Item(159, 4).Value2  = 45180
$ws.Cells.Item(159, 5).Value2  = 13
$ws.Cells.Item(159, 6).Value2  = 100114002
$ws.Cells.Item(159, 7).Value2  = "Camote"
$ws.Cells.Item(159, 8).Value2  = "Sin especificar"
$ws.Cells.Item(159, 9).Value2  = "Primera"
$ws.Cells.Item(159, 10).Value2 = 970
$ws.Cells.Item(159, 11).Value2 = 19000
$ws.Cells.Item(159, 12).Value2 = 20000
$ws.Cells.Item(159, 13).Value2 = 19485
$ws.Cells.Item(159, 14).Value2 = "$/caja 18 kilos"
$ws.Cells.Item(159, 15).Value2 = "Perú"
$ws.Cells.Item(159, 16).Value2 = 1082
$ws.Cells.Item(159, 17).Value2 = 18
$ws.Cells.Item(159, 18).Value2 = "Hortaliza"

# --- New row 160 ---
$ws.Cells.Item(160, 1).Value2  = 9
$ws.Cells.Item(160, 2).Value2  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(160, 3).Value2  = "Metropolitana"
$ws.Cells.Item(160, 4).Value2  = 45180
$ws.Cells.Item(160, 5).Value2  = 13
$ws.Cells.Item(160, 6).Value2  = 100114002
$ws.Cells.Item(160, 7).Value2  = "Camote"
$ws.Cells.Item(160, 8).Value2  = "Sin especificar"
$ws.Cells.Item(160, 9).Value2  = "Primera"
$ws.Cells.Item(160, 10).Value2 = 790
$ws.Cells.Item(160, 11).Value2 = 17000
$ws.Cells.Item(160, 12).Value2 = 18000
$ws.Cells.Item(160, 13).Value2 = 17494
$ws.Cells.Item(160, 14).Value2 = "$/malla 18 kilos"
$ws.Cells.Item(160, 15).Value2 = "Perú"
$ws.Cells.Item(160, 16).Value2 = 972
$ws.Cells.Item(160, 17).Value2 = 18
$ws.Cells.Item(160, 18).Value2 = "Hortaliza"
